$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-02 05:10:17"

$wsZhCn.Range("H2").Value = "2016-09-02 05:10:11"
$wsZhCn.Range("K2").Value = "2016-09-02 05:10:38"

$wsDeDe.Range("K2").Value = "2016-09-02 05:10:44"
